# Update odds values on the active worksheet to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Odd_BTTS_Yes (U) and Odd_BTTS_No (V)
$ws.Range("U6").Value = 1.9
$ws.Range("V6").Value = 1.86

# Row 8: Odd_Over05_FT (M) and Odd_Over15_FT (O)
$ws.Range("M8").Value = 1.03
$ws.Range("O8").Value = 1.22

# Row 9: Odd_Over05_FT (M) and Odd_Over15_FT (O)
$ws.Range("M9").Value = 1.02
$ws.Range("O9").Value = 1.15

# Row 10: Odd_Over05_FT (M) and Odd_Over15_FT (O)
$ws.Range("M10").Value = 1.05
$ws.Range("O10").Value = 1.33
